$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.838.01'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -2.16%  '
$ws.Range("E2").Style = "Normal"

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.238.47'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.30%  '
$ws.Range("E3").Style = "Normal"

# Row 4 - TetherUSD
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("E4").Style = "Normal"

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.78'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.08%  '
$ws.Range("E5").Style = "Normal"

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.98'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -3.75%  '
$ws.Range("E6").Style = "Normal"

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.633'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.58%  '
$ws.Range("E7").Style = "Normal"

# Row 8 - USDC
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.21%  '
$ws.Range("E8").Style = "Normal"

# Row 9 - LidoStakedEther
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.235.92'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.35%  '
$ws.Range("E9").Style = "Normal"

# Row 10 - Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.122'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.95%  '
$ws.Range("E10").Style = "Normal"

# Row 11 - Toncoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.78'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.55%  '
$ws.Range("E11").Style = "Normal"

# Row 12 - Cardano
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.389'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -3.26%  '
$ws.Range("E12").Style = "Normal"

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.800.68'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.34%  '
$ws.Range("E13").Style = "Normal"

# Row 14 - TRON
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.15%  '
$ws.Range("E14").Style = "Normal"

# Row 15 - WrappedBTC
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '64.936.97'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.01%  '
$ws.Range("E15").Style = "Normal"

# Row 16 - Avalanche
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '25.83'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -2.25%  '
$ws.Range("E16").Style = "Normal"

# Row 17 - WrappedEther (was ShibaInu)
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.234.42'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.10%  '
$ws.Range("E17").Style = "Normal"

# Row 18 - ShibaInu (was WrappedEther)
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000159'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.83%  '
$ws.Range("E18").Style = "Normal"

# Row 19 - BitcoinCash
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '416.67'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -4.41%  '
$ws.Range("E19").Style = "Normal"

# Row 20 - Polkadot
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.39'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.34%  '
$ws.Range("E20").Style = "Normal"

# Row 21 - Chainlink
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.83'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.91%  '
$ws.Range("E21").Style = "Normal"

# Row 22 - Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.22'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.49%  '
$ws.Range("E22").Style = "Normal"

# Row 23 - Dai
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("E23").Style = "Normal"

# Row 24 - Litecoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.84'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.38%  '
$ws.Range("E24").Style = "Normal"

# Row 25 - LEO
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.64'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.72%  '
$ws.Range("E25").Style = "Normal"

# Row 26 - Kaspa
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.205'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +4.05%  '
$ws.Range("E26").Style = "Normal"

# Row 27 - Polygon
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.496'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.14%  '
$ws.Range("E27").Style = "Normal"

# Row 28 - PEPE
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.43%  '
$ws.Range("E28").Style = "Normal"

# Row 29 - InternetComputer(DFINITY)
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.07'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.26%  '
$ws.Range("E29").Style = "Normal"

# Row 30 - Binance-PegBSC-USD
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.01%  '
$ws.Range("E30").Style = "Normal"

# Row 31 - PancakeSwap
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -4.12%  '
$ws.Range("E31").Style = "Normal"

# Row 32 - EthereumClassic
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.85'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.08%  '
$ws.Range("E32").Style = "Normal"

# Row 33 - USDe
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.05%  '
$ws.Range("E33").Style = "Normal"

# Row 34 - NEARProtocol
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.00'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.87%  '
$ws.Range("E34").Style = "Normal"

# Row 35 - Aptos
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.43'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.73%  '
$ws.Range("E35").Style = "Normal"

# Row 36 - Fetch.AI
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.51%  '
$ws.Range("E36").Style = "Normal"

# Row 37 - Monero
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '157.49'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.67%  '
$ws.Range("E37").Style = "Normal"

# Row 38 - ImmutableX
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.39'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.08%  '
$ws.Range("E38").Style = "Normal"

# Row 39 - Maker
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.842.70'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.50%  '
$ws.Range("E39").Style = "Normal"

# Row 40 - Stacks
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -3.05%  '
$ws.Range("E40").Style = "Normal"

# Row 41 - EnergySwap
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '25.53'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -4.41%  '
$ws.Range("E41").Style = "Normal"

# Row 42 - Filecoin
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.24'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.41%  '
$ws.Range("E42").Style = "Normal"

# Row 43 - OKB
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '39.59'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.56%  '
$ws.Range("E43").Style = "Normal"

# Row 44 - Mantle
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -6.18%  '
$ws.Range("E44").Style = "Normal"

# Row 45 - RenderToken
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.75'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -4.65%  '
$ws.Range("E45").Style = "Normal"

# Row 46 - Hedera
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0631'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -4.42%  '
$ws.Range("E46").Style = "Normal"

# Row 47 - Bittensor
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '303.23'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -5.72%  '
$ws.Range("E47").Style = "Normal"

# Row 48 - dogwifhat
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -5.28%  '
$ws.Range("E48").Style = "Normal"

# Row 49 - InjectiveProtocol
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '22.13'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -4.93%  '
$ws.Range("E49").Style = "Normal"

# Row 50 - VeChain
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0264'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.12%  '
$ws.Range("E50").Style = "Normal"

# Row 51 - Stellar
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.101'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.23%  '
$ws.Range("E51").Style = "Normal"
